$wb = $excel.ActiveWorkbook

# The shared string "2016-08-19 21:06:48" is used by both the "Overview"
# sheet's "Latest HO Xliff Generate Date" column (G2) and the "de-de"
# sheet's "Correspond Handoff Datetime" column (H2). Update both to keep
# them in sync, matching the single shared-string edit in the diff.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 21:07:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 21:07:30"

# "zh-cn" sheet: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 21:07:25"
$wsZhCn.Range("K2").Value = "2016-08-19 21:07:42"

# "de-de" sheet: Correspond Handback DateTime (K2)
$wsDeDe.Range("K2").Value = "2016-08-19 21:07:49"
